{"js": "// Mill and Ref License Format fix\n// Applies the letter's updated addressee/date/amount details.\n\nconst replacements = [\n  [\"August 05, 2020\", \"August 06, 2020\"],\n  [\"MR. JONATHAN T. GOTIANUN\", \"MS. MINNIE O. CHUA\"],\n  [\"President\", \"President & COO\"],\n  [\"Davao Sugar Central Company, Inc.\", \"Victorias Milling Company, Inc.\"],\n  [\n    \"5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\",\n    \"Ossorio St., Brgy. XVI, Victorias City, Neg. Occ.\",\n  ],\n  [\"Dear Mr. Gotianun:\", \"Dear Pres. Chua:\"],\n  [\n    \"Please be informed that based on your submitted production estimate of 1,000.00 Metric Tons or 1,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is ONE THOUSAND  (PHP 1,000.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of TWO HUNDRED  PESOS (PHP 200.00).\",\n    \"Please be informed that based on your submitted production estimate of 263,250.00 Metric Tons or 5,265,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is FOURTEEN THOUSAND FIVE HUNDRED  (PHP 14,500.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of NINE HUNDRED FIFTY  PESOS (PHP 950.00).\",\n  ],\n  [\n    \"EIGHT HUNDRED  PESOS (PHP 800.00)\",\n    \"THIRTEEN THOUSAND FIVE HUNDRED FIFTY  PESOS (PHP 13,550.00)\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Mill and Ref License Format fix\n# Applies the letter's updated addressee/date/amount details.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"August 05, 2020\"; Replace = \"August 06, 2020\" },\n    @{ Find = \"MR. JONATHAN T. GOTIANUN\"; Replace = \"MS. MINNIE O. CHUA\" },\n    @{ Find = \"President\"; Replace = \"President & COO\" },\n    @{ Find = \"Davao Sugar Central Company, Inc.\"; Replace = \"Victorias Milling Company, Inc.\" },\n    @{ Find = \"5/F, Filinvest Bldg., No. 79 EDSA Highway Hills, Mandaluyong City\"; Replace = \"Ossorio St., Brgy. XVI, Victorias City, Neg. Occ.\" },\n    @{ Find = \"Dear Mr. Gotianun:\"; Replace = \"Dear Pres. Chua:\" },\n    @{ Find = \"Please be informed that based on your submitted production estimate of 1,000.00 Metric Tons or 1,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is ONE THOUSAND  (PHP 1,000.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of TWO HUNDRED  PESOS (PHP 200.00).\"; Replace = \"Please be informed that based on your submitted production estimate of 263,250.00 Metric Tons or 5,265,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is FOURTEEN THOUSAND FIVE HUNDRED  (PHP 14,500.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of NINE HUNDRED FIFTY  PESOS (PHP 950.00).\" },\n    @{ Find = \"EIGHT HUNDRED  PESOS (PHP 800.00)\"; Replace = \"THIRTEEN THOUSAND FIVE HUNDRED FIFTY  PESOS (PHP 13,550.00)\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.Replace\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
